$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 so that all subsequent data rows shift up by one
$ws.Rows.Item(2).Delete()

# Add the new data row at the end (row 21)
$ws.Range("A21").Value = -0.11956262588501
$ws.Range("B21").Value = -0.3192775845527647
$ws.Range("C21").Value = -1.924065947532654
